$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row correct-answer total and the "Total" row values
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 70
$ws.Range("E12").Value = "70/140"
